{"js": "// Replace the multiplication-problem text runs in the table with the\n// new values described by the commit diff. Each old value is unique\n// in the document, so a literal search-and-replace is safe and keeps\n// the original run formatting (font/size) intact.\nconst replacements = [\n  [\"71\u00d712=\", \"68\u00d794=\"],\n  [\"39\u00d736=\", \"36\u00d736=\"],\n  [\"55\u00d715=\", \"14\u00d715=\"],\n  [\"20\u00d753=\", \"93\u00d753=\"],\n  [\"14\u00d718=\", \"62\u00d725=\"],\n  [\"56\u00d760=\", \"37\u00d760=\"],\n  [\"79\u00d711=\", \"58\u00d779=\"],\n  [\"56\u00d791=\", \"95\u00d734=\"],\n  [\"33\u00d726=\", \"12\u00d732=\"],\n  [\"76\u00d749=\", \"87\u00d717=\"],\n  [\"68\u00d750=\", \"81\u00d774=\"],\n  [\"30\u00d750=\", \"95\u00d759=\"],\n  [\"28\u00d769=\", \"21\u00d778=\"],\n  [\"97\u00d754=\", \"33\u00d740=\"],\n  [\"14\u00d747=\", \"83\u00d737=\"],\n  [\"84\u00d715=\", \"33\u00d725=\"],\n  [\"42\u00d771=\", \"94\u00d719=\"],\n  [\"93\u00d723=\", \"56\u00d782=\"],\n  [\"46\u00d796=\", \"95\u00d786=\"],\n  [\"51\u00d743=\", \"73\u00d748=\"],\n  [\"16\u00d791=\", \"71\u00d748=\"],\n  [\"99\u00d770=\", \"74\u00d794=\"],\n  [\"53\u00d757=\", \"26\u00d739=\"],\n  [\"90\u00d758=\", \"54\u00d731=\"],\n  [\"91\u00d720=\", \"31\u00d752=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication-problem text runs in the table with the\n# new values described by the commit diff. Each old value is unique\n# in the document, so a literal Find/Replace (no wildcards) is safe\n# and keeps the original run formatting (font/size) intact.\n$d = $word.ActiveDocument\n\n$olds = @(\n  \"71\u00d712=\",\n  \"39\u00d736=\",\n  \"55\u00d715=\",\n  \"20\u00d753=\",\n  \"14\u00d718=\",\n  \"56\u00d760=\",\n  \"79\u00d711=\",\n  \"56\u00d791=\",\n  \"33\u00d726=\",\n  \"76\u00d749=\",\n  \"68\u00d750=\",\n  \"30\u00d750=\",\n  \"28\u00d769=\",\n  \"97\u00d754=\",\n  \"14\u00d747=\",\n  \"84\u00d715=\",\n  \"42\u00d771=\",\n  \"93\u00d723=\",\n  \"46\u00d796=\",\n  \"51\u00d743=\",\n  \"16\u00d791=\",\n  \"99\u00d770=\",\n  \"53\u00d757=\",\n  \"90\u00d758=\",\n  \"91\u00d720=\"\n)\n$news = @(\n  \"68\u00d794=\",\n  \"36\u00d736=\",\n  \"14\u00d715=\",\n  \"93\u00d753=\",\n  \"62\u00d725=\",\n  \"37\u00d760=\",\n  \"58\u00d779=\",\n  \"95\u00d734=\",\n  \"12\u00d732=\",\n  \"87\u00d717=\",\n  \"81\u00d774=\",\n  \"95\u00d759=\",\n  \"21\u00d778=\",\n  \"33\u00d740=\",\n  \"83\u00d737=\",\n  \"33\u00d725=\",\n  \"94\u00d719=\",\n  \"56\u00d782=\",\n  \"95\u00d786=\",\n  \"73\u00d748=\",\n  \"71\u00d748=\",\n  \"74\u00d794=\",\n  \"26\u00d739=\",\n  \"54\u00d731=\",\n  \"31\u00d752=\"\n)\n\nfor ($i = 0; $i -lt $olds.Count; $i++) {\n  $rng = $d.Content\n  $rng.Find.Execute($olds[$i], $false, $false, $false, $false, $false, $true, 1, $false, $news[$i], 2) | Out-Null\n}\n"}
